# Update start_price sheet
$wb = $excel.ActiveWorkbook

$wsStartPrice = $wb.Worksheets.Item("start_price")
$wsStartPrice.Range("A2").Value = 440.75

$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -3.561601507882178
$wsLinear.Range("B3").Value = 0.0393302453698143
$wsLinear.Range("B4").Value = 896467.2427131189
$wsLinear.Range("B5").Value = "[1.0, 0.14643906734949533, 0.10725670006453873, 0.12915941689363025, 0.12061853767802033, 0.09529222979032297, 0.15027076697280542, 0.28293380902964554, 0.15562662480561312, 0.09925782226734164, 0.10638978403432553, 0.139036908533035, 0.10318990563866105, 0.18032636659089513, 0.23794772145691936, 0.1617859874633006, 0.10038723675731485, 0.14019508733748, 0.10954014483743388, 0.1211875521591939]"

$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 0.8921822897369396
$wsNonLinear.Range("B4").Value = -28.46716274331469
$wsNonLinear.Range("B5").Value = 0.04219520434267009
$wsNonLinear.Range("B6").Value = 1416030.891304129
$wsNonLinear.Range("B7").Value = 38.19857787828995
$wsNonLinear.Range("B8").Value = -0.0829075914580136
$wsNonLinear.Range("B9").Value = 431512.5795807452
$wsNonLinear.Range("B10").Value = "[1.0, 0.14076098385164804, 0.10299558811373223, 0.12411735346320647, 0.11549494245531228, 0.08954633752648712, 0.14580940205569717, 0.2762374586620903, 0.1501245493639528, 0.09487644839189921, 0.10145257250386426, 0.13416056039030638, 0.09915876145921537, 0.17473710792816705, 0.23107156258968406, 0.15577909063882972, 0.09508799457092529, 0.13560128597906634, 0.10414329770282374, 0.11609371422041224]"
